# Update countries & provincias Spain
# Applies the data refresh captured in the authoritative diff:
#  - "Datos actualizados" timestamp bumped from 12:21 to 13:38
#  - Guinea-Bisau's rank moved up (now sits right before Eslovenia),
#    pushing Eslovenia / Islandia / Lituania down one row each
#  - Refreshed totals/actives/recovered/critical/deaths for a handful
#    of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "updated at" timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 13:38"

# --- Row 4: Estados Unidos ---------------------------------------------
$ws.Range("B4").Value = 3696141
$ws.Range("C4").Value = 1116
$ws.Range("D4").Value = 1680424
$ws.Range("E4").Value = 1874587
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 141130

# --- Row 6: India --------------------------------------------------------
$ws.Range("B6").Value = 1009406
$ws.Range("C6").Value = 3769
$ws.Range("D6").Value = 637650
$ws.Range("E6").Value = 346092

# --- Row 19: Alemania ------------------------------------------------------
$ws.Range("B19").Value = 201931
$ws.Range("C19").Value = 95
$ws.Range("E19").Value = 6374

# --- Row 25: Catar -----------------------------------------------------
$ws.Range("B25").Value = 105898
$ws.Range("C25").Value = 421
$ws.Range("D25").Value = 102597
$ws.Range("E25").Value = 3148
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 153

# --- Row 35: Belgica -----------------------------------------------------
$ws.Range("D35").Value = 17269
$ws.Range("E35").Value = 36174

# --- Row 37: Kuwait ------------------------------------------------------
$ws.Range("B37").Value = 58221
$ws.Range("C37").Value = 553
$ws.Range("D37").Value = 48381
$ws.Range("E37").Value = 9436
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 404

# --- Row 53: Suiza -------------------------------------------------------
$ws.Range("B53").Value = 33382
$ws.Range("C53").Value = 92
$ws.Range("E53").Value = 1513

# --- Row 82: Senegal -----------------------------------------------------
$ws.Range("B82").Value = 8544
$ws.Range("C82").Value = 63
$ws.Range("D82").Value = 5809
$ws.Range("E82").Value = 2575
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 160

# --- Row 84: Consejo Danes para los Refugiados ----------------------------
$ws.Range("B84").Value = 8249
$ws.Range("C84").Value = 50
$ws.Range("E84").Value = 3808

# --- Row 88: Finlandia -----------------------------------------------------
$ws.Range("B88").Value = 7301
$ws.Range("C88").Value = 8
$ws.Range("E88").Value = 93

# --- Rows 120-123: Guinea-Bisau overtakes Eslovenia -----------------------
# Row 120 now shows Guinea-Bisau with refreshed figures; Eslovenia,
# Islandia and Lituania each slide down by one row, carrying their
# (former row's) figures with them.
$ws.Range("A120").Value = "Guinea-Bisau"
$ws.Range("B120").Value = 1927
$ws.Range("C120").Value = 25
$ws.Range("D120").Value = 773
$ws.Range("E120").Value = 1128
$ws.Range("H120").Value = 26

$ws.Range("A121").Value = "Eslovenia"
$ws.Range("B121").Value = 1916
$ws.Range("C121").Value = 19
$ws.Range("D121").Value = 1522
$ws.Range("E121").Value = 283
$ws.Range("H121").Value = 111

$ws.Range("A122").Value = "Islandia"
$ws.Range("B122").Value = 1916
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 1895
$ws.Range("E122").Value = 11
$ws.Range("H122").Value = 10

$ws.Range("A123").Value = "Lituania"
$ws.Range("B123").Value = 1908
$ws.Range("C123").Value = 6
$ws.Range("D123").Value = 1595
$ws.Range("E123").Value = 234
$ws.Range("H123").Value = 79

# --- Row 142: Burkina Faso --------------------------------------------------
$ws.Range("B142").Value = 1045
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 887
$ws.Range("E142").Value = 105

# --- Row 147: Surinam ------------------------------------------------------
$ws.Range("B147").Value = 919
$ws.Range("C147").Value = 15
$ws.Range("D147").Value = 582
$ws.Range("E147").Value = 319

# --- Row 155: Malta --------------------------------------------------------
$ws.Range("D155").Value = 662
$ws.Range("E155").Value = 3

# --- Row 162: Vietnam ------------------------------------------------------
$ws.Range("B162").Value = 382
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 26
